# Apply cell updates from the cryptos.xlsx refresh (GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.947.36"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "'1.877.44"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D5").Value = "'0.7444"
$ws.Range("E5").Value = "  -3.61%  "
$ws.Range("D6").Value = "'242.34"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").Value = "'0.9992"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.3162"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'24.86"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07173"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "'0.08444"
$ws.Range("E11").Value = "  -5.37%  "
$ws.Range("D12").Value = "'0.7562"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "'5.440"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "'1.873.53"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "'92.85"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "'29.943.13"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "'6.100"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "'13.65"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "'244.61"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'0.000007842"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "'0.9990"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'2.114.83"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "'8.023"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "'0.9980"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'0.1573"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "'9.335"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").Value = "'164.81"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").Value = "'18.68"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "'2.041"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'1.474"
$ws.Range("E30").Value = "  +3.26%  "
$ws.Range("D31").Value = "'4.615"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").Value = "'1.532"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "'4.285"
$ws.Range("E33").Value = "  +4.42%  "
$ws.Range("D34").Value = "'0.05347"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("D35").Value = "'1.240"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").Value = "'0.7578"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'0.9998"
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "'2.698"
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("D39").Value = "'0.01959"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "'2.750"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "'0.4495"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "'1.110.63"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").Value = "'6.124"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").Value = "'72.79"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").Value = "'0.8640"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'103.32"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "'7.721"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value = "'3.104"
$ws.Range("E49").Value = "  +4.43%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.853"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "'2.013.08"
$ws.Range("E51").Value = "  -0.29%  "
